# Add support for multiple sheets: append a new "Sheet2" after the
# existing "Sheet1", populate it with a 2x2 block of data, and keep
# Sheet1 as the active/selected sheet (matches the source commit).

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Worksheets.Add(Before, After, Count, Type) - insert directly after Sheet1
# so tab order becomes Sheet1, Sheet2 (mirrors Excel's real COM signature).
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "Sheet2"

$sheet2.Range("A1").Value = "a"
$sheet2.Range("B1").Value = "b"
$sheet2.Range("A2").Value = "c"
$sheet2.Range("B2").Value = "d"

# New sheets default to 0.7/0.7/0.75/0.75/0.3/0.3in margins; restore the
# workbook's standard 0.75/0.75/1/1/0.5/0.5in margins (72pt = 1in).
$sheet2.PageSetup.LeftMargin = 54
$sheet2.PageSetup.RightMargin = 54
$sheet2.PageSetup.TopMargin = 72
$sheet2.PageSetup.BottomMargin = 72
$sheet2.PageSetup.HeaderMargin = 36
$sheet2.PageSetup.FooterMargin = 36

# Leave Sheet1 as the selected/active tab, as in the source workbook.
$sheet1.Activate()
